$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at F ("Address"); this pushes the existing
# "District" column (previously F) one column to the right, into G.
$xlShiftToRight = -4161
$ws.Columns("F").Insert($xlShiftToRight)

# New column header
$ws.Range("F2").Value = "Address"

# Populate the new "Address" column (school name/address per teacher row).
# Rows 18, 24, 54 and 56 have no separate address text in the source data,
# so column F is left blank for those rows (same as the un-inserted state).
$ws.Range("F3").Value = 'Govt. High School Badanaguppe'
$ws.Range("F4").Value = 'J S S Girls High School'
$ws.Range("F5").Value = 'G G H S Gundlupet Town'
$ws.Range("F6").Value = 'Govt. High School Gumballi Yelandur'
$ws.Range("F7").Value = 'Govt. High School Haradanahalli'
$ws.Range("F8").Value = 'Govt. High School Kesthur'
$ws.Range("F9").Value = 'Govt. High SchoolNavilur'
$ws.Range("F10").Value = 'Govt. M G S V Jr. CollegeKollegal'
$ws.Range("F11").Value = 'G H S HundipuraGundlupet'
$ws.Range("F12").Value = 'G H S Ginigera'
$ws.Range("F13").Value = 'G H S ArhalGangavathi'
$ws.Range("F14").Value = 'G H S Hiresindogi'
$ws.Range("F15").Value = 'Shree Shashidhara SwamiPU College TavarageraKustagi'
$ws.Range("F16").Value = 'Govt. High School AnegundiGangavathi'
$ws.Range("F17").Value = 'G P U College HanamanalKustagi'
$ws.Range("F19").Value = 'Govt. Higher Primary SchoolSangapurGangavati'
$ws.Range("F20").Value = 'G H S Niralagi'
$ws.Range("F21").Value = 'G H S Hiresindogi'
$ws.Range("F22").Value = 'Govt. High School ChalageraKushtagi'
$ws.Range("F23").Value = 'G H S MandalagiriYelburga'
$ws.Range("F25").Value = 'Govt. H S TalakeriYelburga'
$ws.Range("F26").Value = 'G H S GundagurthiShahapur'
$ws.Range("F27").Value = 'Govt. P U College (H S) B GudiShahapur'
$ws.Range("F28").Value = 'G H S Lingeri Station'
$ws.Range("F29").Value = 'Govt. High School NaganoorShorapur'
$ws.Range("F30").Value = 'G H S Kadechur'
$ws.Range("F31").Value = 'G H P S Filter bed Shahapur'
$ws.Range("F32").Value = 'G H S Yaktapur'
$ws.Range("F33").Value = 'G H S HothpathShahapur'
$ws.Range("F34").Value = 'G H S HalageraShahapur'
$ws.Range("F35").Value = 'G H S TadibidiShahapur'
$ws.Range("F36").Value = 'Maharshi Valmiki GramantarAided High School RukmapurShorapur'
$ws.Range("F37").Value = 'Poojya Shantaveer High SchoolGurmitkal'
$ws.Range("F38").Value = 'G H S BilharShahapur'
$ws.Range("F39").Value = 'G H S YalagiShorapur'
$ws.Range("F40").Value = 'Govt. High School AgniShorapur'
$ws.Range("F41").Value = 'G H S WandurgaShahapur'
$ws.Range("F42").Value = 'G G H S KembhaviShorapur'
$ws.Range("F43").Value = 'G H S Kandkur'
$ws.Range("F44").Value = 'Govt. M P S HalisagarShahapur'
$ws.Range("F45").Value = 'Adarsha VidyalayaShahapura'
$ws.Range("F46").Value = 'Govt. High School (Girls) GogiShahapur'
$ws.Range("F47").Value = 'G H S Basavanthpur'
$ws.Range("F48").Value = 'Govt. High School KolihalShorapur'
$ws.Range("F49").Value = 'Govt. High School WadageraShahapur'
$ws.Range("F50").Value = 'G H S DevargonalShorapur'
$ws.Range("F51").Value = 'Govt. P U College Boys (H S) Shahapur'
$ws.Range("F52").Value = 'G H S JogundabhaviShorapur'
$ws.Range("F53").Value = 'Govt. High School TintaniSurpur'
$ws.Range("F55").Value = 'GHS KannekolurShahapur'
$ws.Range("F57").Value = 'Govt. High SchoolGeddalamariShorapur'
$ws.Range("F58").Value = 'Karanataka Aided High School (UKP) KrishnpurShorapur'
$ws.Range("F59").Value = 'G C P U CollegeKembhaviSurapur'
